# "Tool tip and invalid discount for osprey"
#
# This edit:
#  1. Adds a new column on the "Checkout payments" sheet (inserted before
#     the existing column I), used to carry a new "invalidcode" header /
#     tool-tip style column, with the header label "invalidcode" in I1 and
#     a new "OPSPREY" value in I22 (matching the style of the adjacent H22
#     cell - the Symbol/value pair row).
#  2. Switches the active sheet/selection from "Whishlist" (selecting
#     M11) to "Checkout payments" (selecting H28), which becomes the
#     active tab of the workbook.

$wb = $excel.ActiveWorkbook

# --- Leave the previously-active sheet ("Whishlist") selected on M11 ---
$wsWhish = $wb.Worksheets.Item("Whishlist")
$wsWhish.Activate()
$wsWhish.Range("M11").Select()

# --- Checkout payments: insert a new column before column I ---
$ws = $wb.Worksheets.Item("Checkout payments")
$ws.Columns.Item(9).Insert()

# New header cell for the inserted column
$ws.Range("I1").Value = "invalidcode"

# New value cell in row 22 (same row as the Symbol/value pair near G22/H22)
$ws.Range("I22").Value = "OPSPREY"

# Re-apply the style from the neighboring H22 cell (quotePrefix style),
# since assigning .Value resets the cell style picked up from the
# column insert.
$ws.Range("H22").Copy()
$ws.Range("I22").PasteSpecial(-4122)

# Make "Checkout payments" the active sheet with H28 selected
$ws.Activate()
$ws.Range("H28").Select()
